# Apply weekly crime-data update to the 66th Precinct CompStat workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates -------------------------------------------------
# Volume 31  Number 7 -> Number 8
$ws.Range("A8").Value = "Volume 31   Number  8"

# Report Covering the Week 2/12/2024 Through 2/18/2024 -> 2/19/2024 Through 2/25/2024
$ws.Range("C9").Value = "Report Covering the Week  2/19/2024  Through  2/25/2024"

# --- Row 15: Rape ---------------------------------------------------------
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 500
$ws.Range("N15").Value = 200

# --- Row 16: Robbery ------------------------------------------------------
$ws.Range("C16").Value = "0"
$ws.Range("D16").Value = "0"
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 62.5
$ws.Range("L16").Value = -31.578947368421
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -90.225563909774

# --- Row 17: Fel. Assault --------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 28
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = -6.666666666666
$ws.Range("L17").Value = -20
$ws.Range("M17").Value = 64.705882352941
$ws.Range("N17").Value = -22.222222222222

# --- Row 18: Burglary -------------------------------------------------------
$ws.Range("C18").Value = "0"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 21
$ws.Range("K18").Value = -14.285714285714
$ws.Range("L18").Value = -37.931034482758
$ws.Range("M18").Value = -71.428571428571
$ws.Range("N18").Value = -94.039735099337

# --- Row 19: Gr. Larceny -----------------------------------------------------
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -14.893617021276
$ws.Range("I19").Value = 85
$ws.Range("J19").Value = 94
$ws.Range("K19").Value = -9.574468085106
$ws.Range("L19").Value = 2.409638554216
$ws.Range("M19").Value = 39.344262295082
$ws.Range("N19").Value = -21.296296296296

# --- Row 20: G.L.A. -----------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 6.666666666666
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = 47.619047619047
$ws.Range("L20").Value = 158.333333333333
$ws.Range("M20").Value = 47.619047619047
$ws.Range("N20").Value = -89.491525423728

# --- Row 21: TOTAL --------------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = -6.521739130434
$ws.Range("I21").Value = 183
$ws.Range("J21").Value = 178
$ws.Range("K21").Value = 2.808988764044
$ws.Range("L21").Value = 1.666666666666
$ws.Range("M21").Value = -3.174603174603
$ws.Range("N21").Value = -79.157175398633

# --- Row 22: Transit -----------------------------------------------------------
$ws.Range("G22").Value = "0"
$ws.Range("H22").Value = "***.*"

# --- Row 24: Petit Larceny ------------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -38.095238095238
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = -5.952380952380
$ws.Range("I24").Value = 148
$ws.Range("J24").Value = 172
$ws.Range("K24").Value = -13.953488372093
$ws.Range("L24").Value = -2.631578947368
$ws.Range("M24").Value = -1.333333333333

# --- Row 25: Misd. Assault -------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 73.333333333333
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 39
$ws.Range("K25").Value = 51.282051282051
$ws.Range("L25").Value = 43.902439024390
$ws.Range("M25").Value = -1.666666666666

# --- Row 26: UCR Rape* -------------------------------------------------------------
$ws.Range("C26").Value = "0"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 75
$ws.Range("L26").Value = 75

# --- Row 27: Other Sex Crimes --------------------------------------------------------
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 133.333333333333
$ws.Range("I27").Value = 15
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 50

# --- Row 30: Hate Crimes ----------------------------------------------------------------
$ws.Range("L30").Value = -33.333333333333
